# AFDP-7308 Combine Transcribe and OCR processing into a single media processing module
#
# The "Transcribe Rules" worksheet referenced the old Transcribe model/fact type.
# Update the import of the fact-model class and the rule-table fact declaration
# to point at the new combined MediaEngine model instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Import statement that used to bring in the Transcribe model class now
# imports the combined MediaEngine model class.
$ws.Range("D4").Value = "com.armedia.acm.services.mediaengine.model.MediaEngine"

# Rule-table fact declaration - the $transcribe fact is now typed as MediaEngine.
$ws.Range("C17").Value = "`$transcribe: MediaEngine"
